# Fix de imágenes de los paretos
# Insert a new "Metodo" column at the front of the table on sheet1 (Hoja1),
# label the remaining columns, and fill in the method names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new column before column A; this shifts the existing A:E data to B:F
$ws.Columns.Item(1).Insert()

# Fill column A (new shared strings are registered in this order first)
$ws.Range("A1").Value = "Metodo"
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"

# Updated headers for B1:F1
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

# Column widths matching the real Excel autofit result for the new data
# (values chosen to round-trip as closely as possible through this engine's
# column-width quantization: 12.5 -> 13.33.., 2.3333.. -> 3.1666.., 3.1666.. -> 4)
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(3).ColumnWidth = 3.1666666666666665

$wb.Save()
